{"js": "const pairs = [\n  [\"2023-09-10 Sunday\", \"2023-09-11 Monday\"],\n  [\"34\u00d739=1326\", \"64\u00d789=5696\"],\n  [\"47\u00d760=2820\", \"13\u00d753=689\"],\n  [\"96\u00d799=9504\", \"31\u00d734=1054\"],\n  [\"82\u00d791=7462\", \"59\u00d736=2124\"],\n  [\"97\u00d787=8439\", \"25\u00d731=775\"],\n  [\"51\u00d767=3417\", \"99\u00d716=1584\"],\n  [\"63\u00d762=3906\", \"22\u00d777=1694\"],\n  [\"46\u00d744=2024\", \"44\u00d749=2156\"],\n  [\"61\u00d776=4636\", \"92\u00d766=6072\"],\n  [\"23\u00d744=1012\", \"67\u00d761=4087\"],\n  [\"49\u00d767=3283\", \"34\u00d741=1394\"],\n  [\"71\u00d726=1846\", \"17\u00d723=391\"],\n  [\"51\u00d740=2040\", \"60\u00d729=1740\"],\n  [\"59\u00d790=5310\", \"21\u00d777=1617\"],\n  [\"38\u00d751=1938\", \"70\u00d783=5810\"],\n  [\"72\u00d772=5184\", \"62\u00d763=3906\"],\n  [\"62\u00d736=2232\", \"82\u00d755=4510\"],\n  [\"21\u00d712=252\", \"33\u00d763=2079\"],\n  [\"17\u00d733=561\", \"36\u00d728=1008\"],\n  [\"80\u00d712=960\", \"27\u00d739=1053\"],\n  [\"54\u00d796=5184\", \"73\u00d715=1095\"],\n  [\"80\u00d745=3600\", \"75\u00d733=2475\"],\n  [\"90\u00d722=1980\", \"50\u00d733=1650\"],\n  [\"38\u00d762=2356\", \"52\u00d790=4680\"],\n  [\"23\u00d724=552\", \"15\u00d755=825\"],\n];\n\nconst body = context.document.body;\nfor (const [findText, replaceText] of pairs) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2023-09-10 Sunday\", \"2023-09-11 Monday\"),\n    @(\"34\u00d739=1326\", \"64\u00d789=5696\"),\n    @(\"47\u00d760=2820\", \"13\u00d753=689\"),\n    @(\"96\u00d799=9504\", \"31\u00d734=1054\"),\n    @(\"82\u00d791=7462\", \"59\u00d736=2124\"),\n    @(\"97\u00d787=8439\", \"25\u00d731=775\"),\n    @(\"51\u00d767=3417\", \"99\u00d716=1584\"),\n    @(\"63\u00d762=3906\", \"22\u00d777=1694\"),\n    @(\"46\u00d744=2024\", \"44\u00d749=2156\"),\n    @(\"61\u00d776=4636\", \"92\u00d766=6072\"),\n    @(\"23\u00d744=1012\", \"67\u00d761=4087\"),\n    @(\"49\u00d767=3283\", \"34\u00d741=1394\"),\n    @(\"71\u00d726=1846\", \"17\u00d723=391\"),\n    @(\"51\u00d740=2040\", \"60\u00d729=1740\"),\n    @(\"59\u00d790=5310\", \"21\u00d777=1617\"),\n    @(\"38\u00d751=1938\", \"70\u00d783=5810\"),\n    @(\"72\u00d772=5184\", \"62\u00d763=3906\"),\n    @(\"62\u00d736=2232\", \"82\u00d755=4510\"),\n    @(\"21\u00d712=252\", \"33\u00d763=2079\"),\n    @(\"17\u00d733=561\", \"36\u00d728=1008\"),\n    @(\"80\u00d712=960\", \"27\u00d739=1053\"),\n    @(\"54\u00d796=5184\", \"73\u00d715=1095\"),\n    @(\"80\u00d745=3600\", \"75\u00d733=2475\"),\n    @(\"90\u00d722=1980\", \"50\u00d733=1650\"),\n    @(\"38\u00d762=2356\", \"52\u00d790=4680\"),\n    @(\"23\u00d724=552\", \"15\u00d755=825\"),\n)\n\nforeach ($pair in $pairs) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($find, $false, $true, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
